$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 01:42"
$ws.Range("B4").Value = 85268
$ws.Range("C4").Value = 17057
$ws.Range("E4").Value = 82111
$ws.Range("G4").Value = 266
$ws.Range("H4").Value = 1293
$ws.Range("F21").Value = 23
$ws.Range("F23").Value = 200
$ws.Range("E41").Value = 829
$ws.Range("G41").Value = 5
$ws.Range("H41").Value = 27
$ws.Range("A48").Value = "Panama"
$ws.Range("B48").Value = 674
$ws.Range("C48").Value = 116
$ws.Range("D48").Value = 2
$ws.Range("E48").Value = 663
$ws.Range("F48").Value = 20
$ws.Range("G48").Value = 1
$ws.Range("A49").Value = "Argentina"
$ws.Range("B49").Value = 589
$ws.Range("C49").Value = 87
$ws.Range("D49").Value = 63
$ws.Range("E49").Value = 514
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 4
$ws.Range("H49").Value = 12
$ws.Range("A50").Value = "Peru"
$ws.Range("B50").Value = 580
$ws.Range("C50").Value = 100
$ws.Range("D50").Value = 14
$ws.Range("E50").Value = 557
$ws.Range("F50").Value = 14
$ws.Range("H50").Value = 9
$ws.Range("A51").Value = "Eslovenia"
$ws.Range("B51").Value = 562
$ws.Range("C51").Value = 34
$ws.Range("D51").Value = 10
$ws.Range("E51").Value = 546
$ws.Range("F51").Value = 14
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 6
$ws.Range("A52").Value = "Catar"
$ws.Range("B52").Value = 549
$ws.Range("C52").Value = 12
$ws.Range("D52").Value = 43
$ws.Range("E52").Value = 506
$ws.Range("H52").Value = 0
$ws.Range("A53").Value = "Estonia"
$ws.Range("B53").Value = 538
$ws.Range("C53").Value = 134
$ws.Range("D53").Value = 8
$ws.Range("E53").Value = 529
$ws.Range("F53").Value = 6
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 1
$ws.Range("A74").Value = "Uruguay"
$ws.Range("B74").Value = 238
$ws.Range("C74").Value = 21
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 238
$ws.Range("F74").Value = 3
$ws.Range("H74").Value = 0
$ws.Range("A75").Value = "Costa Rica"
$ws.Range("B75").Value = 231
$ws.Range("C75").Value = 30
$ws.Range("E75").Value = 227
$ws.Range("F75").Value = 5
$ws.Range("H75").Value = 2
$ws.Range("A76").Value = "Eslovaquia"
$ws.Range("B76").Value = 226
$ws.Range("C76").Value = 10
$ws.Range("D76").Value = 2
$ws.Range("E76").Value = 224
$ws.Range("F76").Value = 2
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("A77").Value = "Principado de Andorra"
$ws.Range("B77").Value = 224
$ws.Range("C77").Value = 36
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = 220
$ws.Range("F77").Value = 6
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 3
$ws.Range("A106").Value = "Martinica"
$ws.Range("C106").Value = 15
$ws.Range("E106").Value = 80
$ws.Range("F106").Value = 12
$ws.Range("H106").Value = 1
$ws.Range("A107").Value = "Mauricio"
$ws.Range("B107").Value = 81
$ws.Range("C107").Value = 33
$ws.Range("D107").Value = 0
$ws.Range("E107").Value = 79
$ws.Range("H107").Value = 2
$ws.Range("A108").Value = "Georgia"
$ws.Range("B108").Value = 79
$ws.Range("C108").Value = 4
$ws.Range("D108").Value = 11
$ws.Range("E108").Value = 68
$ws.Range("F108").Value = 1
$ws.Range("A109").Value = "Uzbekistan"
$ws.Range("C109").Value = 15
$ws.Range("D109").Value = 0
$ws.Range("E109").Value = 75
$ws.Range("F109").Value = 4
$ws.Range("H109").Value = 0
$ws.Range("A110").Value = "Camerun"
$ws.Range("B110").Value = 75
$ws.Range("D110").Value = 2
$ws.Range("F110").Value = 0
$ws.Range("A111").Value = "Guadalupe"
$ws.Range("B111").Value = 73
$ws.Range("C111").Value = 0
$ws.Range("E111").Value = 72
$ws.Range("F111").Value = 4
$ws.Range("A112").Value = "Montenegro"
$ws.Range("B112").Value = 69
$ws.Range("C112").Value = 16
$ws.Range("D112").Value = 0
$ws.Range("E112").Value = 68
$ws.Range("F112").Value = 1
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 1
$ws.Range("A113").Value = "Cuba"
$ws.Range("B113").Value = 67
$ws.Range("C113").Value = 10
$ws.Range("D113").Value = 1
$ws.Range("E113").Value = 64
$ws.Range("F113").Value = 2
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 2
$ws.Range("A116").Value = "Bolivia"
$ws.Range("B116").Value = 61
$ws.Range("C116").Value = 29
$ws.Range("E116").Value = 61
$ws.Range("A117").Value = "Liechtenstein"
$ws.Range("B117").Value = 56
$ws.Range("C117").Value = 5
$ws.Range("E117").Value = 56
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 0
$ws.Range("A118").Value = "Honduras"
$ws.Range("B118").Value = 52
$ws.Range("C118").Value = 0
$ws.Range("E118").Value = 51
$ws.Range("H118").Value = 1
$ws.Range("A119").Value = "Consejo Danes para los Refugiados"
$ws.Range("B119").Value = 51
$ws.Range("C119").Value = 3
$ws.Range("E119").Value = 48
$ws.Range("G119").Value = 1
$ws.Range("H119").Value = 3
$ws.Range("A120").Value = "Ruanda"
$ws.Range("B120").Value = 50
$ws.Range("C120").Value = 9
$ws.Range("E120").Value = 50
$ws.Range("A121").Value = "Kirguistan"
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 0
$ws.Range("E121").Value = 44
$ws.Range("F121").Value = 0
$ws.Range("H121").Value = 0
$ws.Range("A122").Value = "Banglades"
$ws.Range("B122").Value = 44
$ws.Range("C122").Value = 5
$ws.Range("D122").Value = 11
$ws.Range("E122").Value = 28
$ws.Range("F122").Value = 1
$ws.Range("H122").Value = 5
$ws.Range("A137").Value = "Barbados"
$ws.Range("C137").Value = 6
$ws.Range("D137").Value = 0
$ws.Range("E137").Value = 24
$ws.Range("A138").Value = "Togo"
$ws.Range("B138").Value = 24
$ws.Range("C138").Value = 1
$ws.Range("D138").Value = 1
$ws.Range("A139").Value = "Madagascar"
$ws.Range("B139").Value = 23
$ws.Range("C139").Value = 4
$ws.Range("E139").Value = 23
$ws.Range("A145").Value = "El Salvador"
$ws.Range("C145").Value = 4
$ws.Range("A146").Value = "Tanzania"
$ws.Range("C146").Value = 0
$ws.Range("A148").Value = "Etiopia"
$ws.Range("C148").Value = 0
$ws.Range("A149").Value = "Guinea Ecuatorial"
$ws.Range("C149").Value = 3
$ws.Range("A150").Value = "Republica de Yibuti"
$ws.Range("A151").Value = "Mongolia"
$ws.Range("C151").Value = 1
$ws.Range("A152").Value = "San Martin (Parte Francesa)"
$ws.Range("A153").Value = "Dominica"
$ws.Range("C153").Value = 0
$ws.Range("A156").Value = "Haiti"
$ws.Range("A157").Value = "Surinam"
$ws.Range("A160").Value = "Antigua y Barbuda"
$ws.Range("C160").Value = 4
$ws.Range("A161").Value = "Mozambique"
$ws.Range("C161").Value = 2
$ws.Range("A162").Value = "Granada"
$ws.Range("C162").Value = 6
$ws.Range("A163").Value = "Seychelles"
$ws.Range("C163").Value = 0
$ws.Range("A165").Value = "Curazao"
$ws.Range("B165").Value = 7
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 2
$ws.Range("E165").Value = 4
$ws.Range("H165").Value = 1
$ws.Range("A166").Value = "Eritrea"
$ws.Range("C166").Value = 2
$ws.Range("A167").Value = "Benin"
$ws.Range("C167").Value = 0
$ws.Range("A168").Value = "Suazilandia"
$ws.Range("C168").Value = 2
$ws.Range("A169").Value = "Laos"
$ws.Range("C169").Value = 3
$ws.Range("D169").Value = 0
$ws.Range("E169").Value = 6
$ws.Range("A170").Value = "Groenlandia"
$ws.Range("E170").Value = 4
$ws.Range("H170").Value = 0
$ws.Range("A171").Value = "Siria"
$ws.Range("C171").Value = 0
$ws.Range("A172").Value = "Birmania"
$ws.Range("C172").Value = 2
$ws.Range("A175").Value = "Cabo Verde"
$ws.Range("C175").Value = 1
$ws.Range("A176").Value = "Guyana"
$ws.Range("B176").Value = 5
$ws.Range("H176").Value = 1
$ws.Range("A177").Value = "Guinea"
$ws.Range("C177").Value = 0
$ws.Range("A178").Value = "Angola"
$ws.Range("C178").Value = 1
$ws.Range("A179").Value = "Congo"
$ws.Range("C179").Value = 0
$ws.Range("A181").Value = "Mali"
$ws.Range("C181").Value = 2
$ws.Range("E181").Value = 4
$ws.Range("H181").Value = 0
$ws.Range("A182").Value = "Republica de Africa Central"
$ws.Range("A183").Value = "Mauritania"
$ws.Range("C183").Value = 1
$ws.Range("A184").Value = "San Bartolome"
$ws.Range("A185").Value = "San Martin (Parte Holandesa)"
$ws.Range("A186").Value = "Liberia"
$ws.Range("A187").Value = "Republica del Chad"
$ws.Range("C187").Value = 0
$ws.Range("A188").Value = "Gambia"
$ws.Range("E188").Value = 2
$ws.Range("H188").Value = 1
$ws.Range("A189").Value = "Santa Lucia"
$ws.Range("D189").Value = 1
$ws.Range("H189").Value = 0
$ws.Range("A192").Value = "Sudan"
$ws.Range("A203").Value = "Timor Oriental"
$ws.Range("A204").Value = "San Vicente y las Granadinas"
